# Update cryptos list with refreshed prices/volumes scraped on Fri Jul 28 13:26:43 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{ Row=2; D='29.280.82' },
    @{ Row=3; D='1.870.07'; E='  -0.52%  ' },
    @{ Row=4; E='  -0.22%  ' },
    @{ Row=5; D='0.7115'; E='  -0.86%  ' },
    @{ Row=6; D='241.76'; E='  -0.06%  ' },
    @{ Row=7; E='  -0.20%  ' },
    @{ Row=8; D='0.3107'; E='  +0.24%  ' },
    @{ Row=9; D='0.07747'; E='  -2.54%  ' },
    @{ Row=10; D='24.75'; E='  -2.32%  ' },
    @{ Row=11; E='  +1.56%  ' },
    @{ Row=12; D='1.894.15'; E='  -0.01%  ' },
    @{ Row=13; D='5.222'; E='  -1.15%  ' },
    @{ Row=14; D='0.7117'; E='  -2.44%  ' },
    @{ Row=15; D='91.07'; E='  -0.11%  ' },
    @{ Row=16; D='29.289.93'; E='  -0.83%  ' },
    @{ Row=17; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000008230'; E='  +4.46%  ' },
    @{ Row=18; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='5.938'; E='  +0.52%  ' },
    @{ Row=19; D='243.59'; E='  -1.00%  ' },
    @{ Row=20; D='2.122.31'; E='  -0.61%  ' },
    @{ Row=21; D='13.15'; E='  -1.38%  ' },
    @{ Row=22; D='0.9998'; E='  -0.30%  ' },
    @{ Row=23; D='7.865'; E='  -2.50%  ' },
    @{ Row=24; E='  -0.21%  ' },
    @{ Row=25; D='0.1631'; E='  +1.42%  ' },
    @{ Row=26; D='163.71'; E='  -0.04%  ' },
    @{ Row=27; D='9.002'; E='  -0.50%  ' },
    @{ Row=28; D='18.52'; E='  +1.02%  ' },
    @{ Row=29; D='1.509'; E='  +0.97%  ' },
    @{ Row=30; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='1.310'; E='  -3.44%  ' },
    @{ Row=31; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='4.400'; E='  -0.02%  ' },
    @{ Row=32; D='4.272'; E='  +3.99%  ' },
    @{ Row=33; D='0.05168'; E='  -0.81%  ' },
    @{ Row=34; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.920'; E='  -1.32%  ' },
    @{ Row=35; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.7753'; E='  +6.75%  ' },
    @{ Row=36; D='1.171'; E='  -2.35%  ' },
    @{ Row=37; D='2.688'; E='  +0.11%  ' },
    @{ Row=38; D='0.01861'; E='  -0.58%  ' },
    @{ Row=39; D='2.712'; E='  +0.11%  ' },
    @{ Row=40; D='1.161.57'; E='  -3.61%  ' },
    @{ Row=41; D='6.411'; E='  +3.93%  ' },
    @{ Row=42; D='73.32'; E='  -0.40%  ' },
    @{ Row=43; D='0.8907'; E='  -2.20%  ' },
    @{ Row=44; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='104.86'; E='  +2.45%  ' },
    @{ Row=45; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='1.000'; E='  -0.25%  ' },
    @{ Row=46; D='2.020.21'; E='  +0.07%  ' },
    @{ Row=47; D='1.793'; E='  -0.22%  ' },
    @{ Row=48; D='0.5189'; E='  -2.08%  ' },
    @{ Row=49; D='9.383'; E='  +0.75%  ' },
    @{ Row=50; D='0.00000000119'; E='  +1.98%  ' },
    @{ Row=51; D='0.4298'; E='  -0.65%  ' }
)

# Columns D (Price) and E (Volume) hold number-shaped text (e.g. "29.280.82",
# "1.310", "0.000008230") that must stay text, not be coerced to numbers -
# so those two columns are written with a leading apostrophe, just as typing
# them in the Excel UI would.
$textColumns = @("D", "E")

foreach ($update in $rowUpdates) {
    $row = $update.Row
    foreach ($col in @("A", "B", "C", "D", "E")) {
        if ($update.ContainsKey($col)) {
            $value = $update[$col]
            $cell = $ws.Range("$col$row")
            if ($textColumns -contains $col) {
                $cell.Value = "'" + $value
            } else {
                $cell.Value = $value
            }
        }
    }
}

